# Update TestExcel.xlsx: replace the numeric Id column with GUID-style string
# identifiers, widen column A to fit them, set the page setup for printing,
# and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric Id values in A2:A4 with GUID strings.
$ws.Range("A2").Value = "8BC78143-9FD5-45E4-AEED-F5648D58473C"
$ws.Range("A3").Value = "46C6F115-B719-48BF-8EE1-3ABF480DF748"
$ws.Range("A4").Value = "5088AB6B-CFCE-4531-BDFE-1E79CCAA7A3D"

# Widen column A so the GUIDs are fully visible.
$ws.Columns.Item(1).ColumnWidth = 39

# Configure the page setup for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the selection to match where the author left off editing.
$ws.Range("B10").Select() | Out-Null
